$d = $word.ActiveDocument

# 1. Consolidate the title, author and abstract text into single runs
#    (the underlying multi-run spellings collapse to one run each).
$d.Content.Find.Execute("Questions: Arithmetic on complex numbers", $false, $false, $false, $false, $false, $true, 1, $false, "Questions: Arithmetic on complex numbers", 2) | Out-Null
$d.Content.Find.Execute("Charlotte McCarthy", $false, $false, $false, $false, $false, $true, 1, $false, "Charlotte McCarthy", 2) | Out-Null
$d.Content.Find.Execute("A selection of questions for the study guide on arithmetic on complex numbers.", $false, $false, $false, $false, $false, $true, 1, $false, "A selection of questions for the study guide on arithmetic on complex numbers.", 2) | Out-Null

# 2. Re-order the <m:sepChr/> element to sit between <m:begChr/> and
#    <m:endChr/> inside every bracket (m:d) delimiter's properties
#    (m:dPr), for every math region in the document.
$oldPr = '<m:dPr><m:begChr m:val="(" /><m:endChr m:val=")" /><m:sepChr m:val="" /><m:grow /></m:dPr>'
$newPr = '<m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr>'

$count = $d.OMaths.Count
for ($i = 1; $i -le $count; $i++) {
    $om = $d.OMaths.Item($i)
    $xml = $om.Range.WordOpenXML
    $start = $xml.IndexOf("<m:oMath>")
    $end = $xml.IndexOf("</m:oMath>")
    if ($start -lt 0 -or $end -lt 0) {
        continue
    }
    $frag = $xml.Substring($start, $end - $start + 10)
    if ($frag.Contains("m:dPr")) {
        $newFrag = $frag.Replace($oldPr, $newPr)
        $om.Range.InsertXML($newFrag)
    }
}
